# Natmi following Dr Hou advice
# Recomputes the Thbs1 -> Tnfrsf11b sending/target cluster table for ECs, FAPs, sCs
# (previously only the matching sending/target cluster rows existed; now the full
# 3x3 cross of sending x target clusters is present, with updated statistics).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Sending cluster, Target cluster, then E..T numeric values
$rows = @(
    @{ Row=2;  A="ECs";  D="ECs";  E=3; F=1; G=142.9073533333333; H=428.72206;     I=0.5576664151504187; J=0.5576664151504188; K=1; L=0.3333333333333333; M=0.1176943333333333; N=0.353083;  O=0.05829606481791055; P=0.05829606481791055; Q=16.81938567899778;  R=151.37447111098;    S=0.03250975748438062; T=0.03250975748438063 }
    @{ Row=3;  A="ECs";  D="FAPs"; E=3; F=1; G=142.9073533333333; H=428.72206;     I=0.5576664151504187; J=0.5576664151504188; K=3; L=1;                  M=1.183046666666667;  N=3.54914;   O=0.5859837360842607;  P=0.5859837360842608;  Q=169.0660680031556;  R=1521.5946120284;    S=0.3267834494385587;  T=0.3267834494385589 }
    @{ Row=4;  A="ECs";  D="sCs";  E=3; F=1; G=142.9073533333333; H=428.72206;     I=0.5576664151504187; J=0.5576664151504188; K=3; L=1;                  M=0.718166;           N=2.154498;  O=0.3557201990978286;  P=0.3557201990978286;  Q=102.6312023139867;  R=923.6808208258799;  S=0.1983732082274793;  T=0.1983732082274793 }
    @{ Row=5;  A="FAPs"; D="ECs";  E=3; F=1; G=63.967809;         H=191.903427;    I=0.2496211559306514; J=0.2496211559306514; K=1; L=0.3333333333333333; M=0.1176943333333333; N=0.353083;  O=0.05829606481791055; P=0.05829606481791055; Q=7.528648635049;     R=67.757837715441;    S=0.01455193108605501; T=0.01455193108605501 }
    @{ Row=6;  A="FAPs"; D="FAPs"; E=3; F=1; G=63.967809;         H=191.903427;    I=0.2496211559306514; J=0.2496211559306514; K=3; L=1;                  M=1.183046666666667;  N=3.54914;   O=0.5859837360842607;  P=0.5859837360842608;  Q=75.67690321141998;  R=681.09212890278;    S=0.1462739375579149;  T=0.1462739375579149 }
    @{ Row=7;  A="FAPs"; D="sCs";  E=3; F=1; G=63.967809;         H=191.903427;    I=0.2496211559306514; J=0.2496211559306514; K=3; L=1;                  M=0.718166;           N=2.154498;  O=0.3557201990978286;  P=0.3557201990978286;  Q=45.939505518294;    R=413.455549664646;   S=0.08879528728668144; T=0.08879528728668144 }
    @{ Row=8;  A="sCs";  D="ECs";  E=3; F=1; G=49.38440333333333; H=148.15321;     I=0.1927124289189298; J=0.1927124289189298; K=1; L=0.3333333333333333; M=0.1176943333333333; N=0.353083;  O=0.05829606481791055; P=0.05829606481791055; Q=5.812264427381111;  R=52.31037984643;     S=0.01123437624747491; T=0.01123437624747491 }
    @{ Row=9;  A="sCs";  D="FAPs"; E=3; F=1; G=49.38440333333333; H=148.15321;     I=0.1927124289189298; J=0.1927124289189298; K=3; L=1;                  M=1.183046666666667;  N=3.54914;   O=0.5859837360842607;  P=0.5859837360842608;  Q=58.42405374882222;  R=525.8164837394;     S=0.112926349087787;   T=0.1129263490877871 }
    @{ Row=10; A="sCs";  D="sCs";  E=3; F=1; G=49.38440333333333; H=148.15321;     I=0.1927124289189298; J=0.1927124289189298; K=3; L=1;                  M=0.718166;           N=2.154498;  O=0.3557201990978286;  P=0.3557201990978286;  Q=35.46619940428666;  R=319.19579463858;    S=0.06855170358366787; T=0.06855170358366787 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = "Thbs1"
    $ws.Range("C$row").Value = "Tnfrsf11b"
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
}
